# excel data source files udpates
#
# Adds two new data-source sheets ("weather", "altitude") after the
# existing "lapsundercut" sheet, and adds two new columns
# (tyre_before / tyre_after) to the "pitstop" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) pitstop: add tyre_before / tyre_after columns (H, I)
# ---------------------------------------------------------------
$pitstop = $wb.Worksheets.Item("pitstop")

$pitstop.Cells.Item(1, 8).Value = "tyre_before"
$pitstop.Cells.Item(1, 9).Value = "tyre_after"

$tyreData = @(
    @(4,3), @(3,2), @(4,2), @(2,2), @(2,3), @(3,3), @(3,4), @(4,3), @(3,4), @(3,4),
    @(3,4), @(4,2), @(4,4), @(4,3), @(4,3), @(3,4), @(4,3), @(3,4), @(4,3), @(3,4),
    @(4,3), @(4,3), @(3,4), @(4,3), @(3,2), @(4,3), @(3,4), @(4,3), @(4,3)
)
for ($i = 0; $i -lt $tyreData.Length; $i++) {
    $row = $i + 2
    $pair = $tyreData[$i]
    $pitstop.Cells.Item($row, 8).Value = $pair[0]
    $pitstop.Cells.Item($row, 9).Value = $pair[1]
}

$pitstop.Columns.Item(8).ColumnWidth = 10.666666666666666
$pitstop.Columns.Item(9).ColumnWidth = 9

$pitstop.Range("H31").Select()

# ---------------------------------------------------------------
# 2) new "weather" sheet, placed after the last existing sheet
# ---------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$weather = $wb.Worksheets.Add([Type]::Missing, $lastSheet)
$weather.Name = "weather"

$weather.Range("A1").Value = "Skycondition"
$weather.Range("B1").Value = "Partly Cloudy"
$weather.Range("A2").Value = "Precipation type"
$weather.Range("B2").Value = "rain"
$weather.Range("A3").Value = "Temperature"
$weather.Range("B3").Value = "66.84°F"
$weather.Range("A4").Value = "Humidity"
$weather.Range("B4").Value = 0.57
$weather.Range("B4").NumberFormat = "0%"
$weather.Range("A5").Value = "Wind speed"
$weather.Range("B5").Value = "7.48 mph"
$weather.Range("A6").Value = "Wind bearing"
$weather.Range("B6").Value = "273°"

$weather.Columns.Item(1).ColumnWidth = 14.833333333333334
$weather.Columns.Item(2).ColumnWidth = 11.833333333333334

$weather.PageSetup.Orientation = 1

$weather.Range("C2").Select()

# ---------------------------------------------------------------
# 3) new "altitude" sheet, placed after "weather"
# ---------------------------------------------------------------
$altitude = $wb.Worksheets.Add([Type]::Missing, $weather)
$altitude.Name = "altitude"

$altitude.Range("A1").Value = "delta"
$altitude.Range("B1").Value = 12.8

$altitude.Range("B2").Select()
